$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '26.384.37'
$ws.Range('E2').Value2 = '  -1.99%  '

$ws.Range('D3').Value2 = '1.838.94'
$ws.Range('E3').Value2 = '  -2.28%  '

$ws.Range('D4').Value2 = "'1.002"

$ws.Range('D5').Value2 = "'259.42"
$ws.Range('E5').Value2 = '  -6.85%  '

$ws.Range('E6').Value2 = '  +0.20%  '

$ws.Range('D7').Value2 = "'0.5212"
$ws.Range('E7').Value2 = '  -2.91%  '

$ws.Range('D8').Value2 = "'0.3240"
$ws.Range('E8').Value2 = '  -6.25%  '

$ws.Range('D9').Value2 = "'0.06757"
$ws.Range('E9').Value2 = '  -3.22%  '

$ws.Range('D10').Value2 = "'18.47"
$ws.Range('E10').Value2 = '  -8.45%  '

$ws.Range('D11').Value2 = "'0.7616"
$ws.Range('E11').Value2 = '  -5.82%  '

$ws.Range('D12').Value2 = "'0.07678"
$ws.Range('E12').Value2 = '  -0.43%  '

$ws.Range('D13').Value2 = '1.813.58'
$ws.Range('E13').Value2 = '  -3.53%  '

$ws.Range('D14').Value2 = "'87.90"
$ws.Range('E14').Value2 = '  -3.40%  '

$ws.Range('D15').Value2 = "'5.006"
$ws.Range('E15').Value2 = '  -3.51%  '

$ws.Range('D16').Value2 = "'1.003"
$ws.Range('E16').Value2 = '  +0.34%  '

$ws.Range('D17').Value2 = "'13.85"
$ws.Range('E17').Value2 = '  -5.11%  '

$ws.Range('E18').Value2 = '  +0.23%  '

$ws.Range('D19').Value2 = "'0.000007945"
$ws.Range('E19').Value2 = '  -1.24%  '

$ws.Range('D20').Value2 = '26.464.52'
$ws.Range('E20').Value2 = '  -1.89%  '

$ws.Range('D21').Value2 = '2.080.23'
$ws.Range('E21').Value2 = '  -1.37%  '

$ws.Range('D22').Value2 = "'4.551"
$ws.Range('E22').Value2 = '  -4.29%  '

$ws.Range('D23').Value2 = "'9.431"
$ws.Range('E23').Value2 = '  -6.43%  '

$ws.Range('D24').Value2 = "'5.928"
$ws.Range('E24').Value2 = '  -4.39%  '

$ws.Range('D25').Value2 = "'144.01"
$ws.Range('E25').Value2 = '  -2.04%  '

$ws.Range('D26').Value2 = "'2.196"
$ws.Range('E26').Value2 = '  -7.81%  '

$ws.Range('D27').Value2 = "'1.643"
$ws.Range('E27').Value2 = '  -1.30%  '

$ws.Range('D28').Value2 = "'16.94"
$ws.Range('E28').Value2 = '  -2.54%  '

$ws.Range('D29').Value2 = "'111.08"
$ws.Range('E29').Value2 = '  -2.52%  '

$ws.Range('D30').Value2 = "'4.140"
$ws.Range('E30').Value2 = '  -5.31%  '

$ws.Range('D31').Value2 = "'4.108"
$ws.Range('E31').Value2 = '  -5.11%  '

$ws.Range('D32').Value2 = "'0.08711"
$ws.Range('E32').Value2 = '  -2.42%  '

$ws.Range('D33').Value2 = "'0.04756"
$ws.Range('E33').Value2 = '  -3.83%  '

$ws.Range('D34').Value2 = "'1.118"
$ws.Range('E34').Value2 = '  -5.38%  '

$ws.Range('D35').Value2 = "'2.849"
$ws.Range('E35').Value2 = '  -1.23%  '

$ws.Range('D36').Value2 = "'0.6953"
$ws.Range('E36').Value2 = '  -5.52%  '

$ws.Range('D37').Value2 = "'3.060"
$ws.Range('E37').Value2 = '  -7.15%  '

$ws.Range('D38').Value2 = "'0.01747"
$ws.Range('E38').Value2 = '  -5.74%  '

$ws.Range('D39').Value2 = "'2.153"
$ws.Range('E39').Value2 = '  -9.58%  '

$ws.Range('D40').Value2 = "'0.4807"
$ws.Range('E40').Value2 = '  -6.99%  '

$ws.Range('D41').Value2 = "'111.05"
$ws.Range('E41').Value2 = '  -4.17%  '

$ws.Range('D42').Value2 = "'0.8887"
$ws.Range('E42').Value2 = '  -7.38%  '

$ws.Range('D43').Value2 = "'6.049"
$ws.Range('E43').Value2 = '  -2.44%  '

$ws.Range('E44').Value2 = '  +0.30%  '

$ws.Range('D45').Value2 = "'7.607"
$ws.Range('E45').Value2 = '  -6.65%  '

$ws.Range('B46').Value2 = 'Cronos'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value2 = "'0.05856"
$ws.Range('E46').Value2 = '  -1.66%  '

$ws.Range('B47').Value2 = 'Decentraland'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value2 = "'0.4107"
$ws.Range('E47').Value2 = '  -8.83%  '

$ws.Range('D48').Value2 = "'8.933"
$ws.Range('E48').Value2 = '  -4.55%  '

$ws.Range('D49').Value2 = "'34.94"
$ws.Range('E49').Value2 = '  -3.92%  '

$ws.Range('D50').Value2 = "'0.1214"
$ws.Range('E50').Value2 = '  -10.00%  '

$ws.Range('D51').Value2 = "'0.8825"
$ws.Range('E51').Value2 = '  -0.48%  '
